# Commit: "Added data and downloader"
# The factor labels in column A of the "modENCODE data" sheet are prefixed
# with "#" so a downstream downloader script can recognize/flag them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("modENCODE data")

# ELT-2 -> #ELT-2
$ws.Range("A2").Value = "#ELT-2"
$ws.Range("A4").Value = "#ELT-2"
$ws.Range("A6").Value = "#ELT-2"

# DVE-1 -> #DVE-1
$ws.Range("A8").Value = "#DVE-1"
$ws.Range("A9").Value = "#DVE-1"
$ws.Range("A10").Value = "#DVE-1"
$ws.Range("A11").Value = "#DVE-1"

# PHA-4 -> #PHA-4
$ws.Range("A12").Value = "#PHA-4"
$ws.Range("A13").Value = "#PHA-4"
$ws.Range("A14").Value = "#PHA-4"
$ws.Range("A15").Value = "#PHA-4"

# SMA-9 -> #SMA-9
$ws.Range("A18").Value = "#SMA-9"

# Egl-27 -> #Egl-27
$ws.Range("A19").Value = "#Egl-27"

# NHR-2 -> #NHR-2
$ws.Range("A20").Value = "#NHR-2"
$ws.Range("A21").Value = "#NHR-2"
$ws.Range("A22").Value = "#NHR-2"
$ws.Range("A23").Value = "#NHR-2"

# Update the active cell/selection to match the author's saved cursor position
$ws.Range("A10").Select()
